# Insert a new "Required Software" section before the existing
# "Required Text" heading, per the commit "added software install to cs
# syllabus".
#
# Strategy: locate the "Required Text" Heading 1 paragraph, insert all of
# the new paragraphs' plain text (with embedded carriage returns) right
# before it in one shot (so the paragraph count/order come out correct and
# nothing bleeds into neighboring paragraphs), then walk back over the
# newly created paragraphs to apply the correct styles, run formatting,
# and hyperlinks.

$d = $word.ActiveDocument
$cr = [char]13

# --- find the "Required Text" heading paragraph -----------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd($cr) -eq "Required Text" -and $p.Style.NameLocal -eq "Heading 1") {
        $target = $p
        break
    }
}

$insertStart = $target.Range.Start

# --- build the plain-text block (one line per new paragraph) ----------
$lines = @(
    "Required Software",
    "For this class we will be programming in the Python programming language, using a development platform called Jupyter Notebook.",
    "In order to run this software, it is strongly recommended that you use a computer running Windows, MacOS, or a desktop Linux (e.g. Ubuntu, Debian). If you are working from an iPad, Chromebook, or Android the configuration will be more difficult.",
    "The easiest way to install the latest version of Python 3 and Jupyter is by downloading and running the graphical installer for Anaconda. Anaconda is a complete data science platform, but it contains everything we need in a neat package.",
    "**Click here to find the Anaconda installer for your platform."
)

$block = ($lines -join $cr) + $cr

$ins = $d.Range($insertStart, $insertStart)
$ins.InsertBefore($block)

# --- figure out which paragraph indices the new block landed on -------
$firstIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Start -eq $insertStart) {
        $firstIndex = $i
        break
    }
}

$pHeading      = $d.Paragraphs($firstIndex)
$pIntro        = $d.Paragraphs($firstIndex + 1)
$pRecommend    = $d.Paragraphs($firstIndex + 2)
$pInstall      = $d.Paragraphs($firstIndex + 3)
$pInstallLink  = $d.Paragraphs($firstIndex + 4)

# --- paragraph styles ---------------------------------------------------
$pHeading.Style     = "Heading 1"
$pIntro.Style       = "First Paragraph"
$pRecommend.Style   = "Body Text"
$pInstall.Style     = "Body Text"
$pInstallLink.Style = "Body Text"

# --- bookmark for the new heading ---------------------------------------
# (pass a standalone Range, not the paragraph's live .Range, or the
# bookmark collection re-sorts/recomputes against a stale anchor)
$headingBmRange = $d.Range($pHeading.Range.Start, $pHeading.Range.End)
$d.Bookmarks.Add("required-software", $headingBmRange) | Out-Null

# --- hyperlink: "Python" -> https://python.org (in the intro paragraph) -
$introRange = $pIntro.Range
$introRange.Find.ClearFormatting()
$introRange.Find.Execute("Python", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Hyperlinks.Add($introRange, "https://python.org", "", "", "Python") | Out-Null

# --- hyperlink: "Jupyter Notebook" -> https://jupyter.org/ --------------
$jnRange = $d.Range(0, $d.Content.End)
$jnRange.Find.ClearFormatting()
$jnRange.Find.Execute("Jupyter Notebook", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Hyperlinks.Add($jnRange, "https://jupyter.org/", "", "", "Jupyter Notebook") | Out-Null

# --- italicize "Anaconda" (first occurrence, in the install paragraph) --
$anRange = $d.Range(0, $d.Content.End)
$anRange.Find.ClearFormatting()
$anRange.Find.Execute("Anaconda", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anRange.Italic = 1

# --- hyperlink the "Click here..." sentence in the last new paragraph ---
$linkText = "Click here to find the Anaconda installer for your platform."
$clickRange = $d.Range(0, $d.Content.End)
$clickRange.Find.ClearFormatting()
$clickRange.Find.Execute($linkText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Hyperlinks.Add($clickRange, "https://www.anaconda.com/download", "", "", $linkText) | Out-Null
